# Update weekly fruit/vegetable price data (Fruta, Terminal Hortofrutícola Agro Chillán - Damasco)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44552
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 15500
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15750
$ws.Range("S2").Value = 1050

# Row 3
$ws.Range("D3").Value = 44537
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 21500
$ws.Range("P3").Value = 21250
$ws.Range("S3").Value = 1417

# Row 4
$ws.Range("D4").Value = 44176
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17400
$ws.Range("Q4").Value = "`$/caja 18 kilos"
$ws.Range("S4").Value = 967
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44165
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 17000
$ws.Range("P5").Value = 16500
$ws.Range("Q5").Value = "`$/caja 15 kilos granel"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 1100
$ws.Range("T5").Value = 15

# Row 6
$ws.Range("D6").Value = 44907
$ws.Range("K6").Value = "Castle Brite"
$ws.Range("M6").Value = 120
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = "`$/bandeja 10 kilos"
$ws.Range("S6").Value = 1550
$ws.Range("T6").Value = 10

# Row 7
$ws.Range("D7").Value = 44907
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 14000
$ws.Range("Q7").Value = "`$/bandeja 10 kilos"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1400
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("D8").Value = 44187
$ws.Range("M8").Value = 55
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15455
$ws.Range("Q8").Value = "`$/caja 15 kilos granel"
$ws.Range("S8").Value = 1030
$ws.Range("T8").Value = 15

# Row 9
$ws.Range("D9").Value = 44168
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 16000
$ws.Range("O9").Value = 17000
$ws.Range("P9").Value = 16500
$ws.Range("Q9").Value = "`$/caja 16 kilos granel"
$ws.Range("R9").Value = "Región de Coquimbo"
$ws.Range("S9").Value = 1031
$ws.Range("T9").Value = 16

# Row 10
$ws.Range("D10").Value = 44544
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 17000
$ws.Range("P10").Value = 16500
$ws.Range("S10").Value = 1100

# Row 11
$ws.Range("D11").Value = 44174
$ws.Range("M11").Value = 75
$ws.Range("N11").Value = 9000
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 9467
$ws.Range("Q11").Value = "`$/caja 10 kilos"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 947
$ws.Range("T11").Value = 10

# Row 12
$ws.Range("D12").Value = 44551
$ws.Range("M12").Value = 120
$ws.Range("N12").Value = 15500
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 15750
$ws.Range("S12").Value = 1050

# Row 13
$ws.Range("D13").Value = 44181
$ws.Range("K13").Value = "Modesto"
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 21000
$ws.Range("P13").Value = 20500
$ws.Range("Q13").Value = "`$/caja 18 kilos"
$ws.Range("R13").Value = "Región de Coquimbo"
$ws.Range("S13").Value = 1139
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("D14").Value = 44904
$ws.Range("L14").Value = "Primera"
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 16000
$ws.Range("P14").Value = 15500
$ws.Range("S14").Value = 1550

# Row 15
$ws.Range("D15").Value = 44904
$ws.Range("K15").Value = "Castle Brite"
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 14000
$ws.Range("P15").Value = 14000
$ws.Range("Q15").Value = "`$/bandeja 10 kilos"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 1400
$ws.Range("T15").Value = 10

# Row 16
$ws.Range("D16").Value = 44189
$ws.Range("K16").Value = "Dina"
$ws.Range("M16").Value = 80
$ws.Range("N16").Value = 16000
$ws.Range("O16").Value = 17000
$ws.Range("P16").Value = 16562
$ws.Range("S16").Value = 920
